# FAST_holdings.xlsx update:
#  - refresh the "as of" date in the confidential disclaimer text (A13)
#  - refresh the Weight (D) / Percent Change (E) values for rows 2-10
#
# The sheet ships protected (legacy password hash), so unlock it first,
# make the data edits, and leave the data itself untouched everywhere else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Disclaimer text: bump the "as of" date from 2021-03-19 to 2021-03-22
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-22 for illustrative purposes only and are subject to change."

# Row 2 (ARKG)
$ws.Range("D2").Value = 0.1085888285438651
$ws.Range("E2").Value = 0.02255970067128854

# Row 3 (MTUM)
$ws.Range("D3").Value = 0.1040990287056252
$ws.Range("E3").Value = 0.01412447190869548

# Row 4 (QUAL)
$ws.Range("D4").Value = 0.1148448586828432
$ws.Range("E4").Value = 0.008938267479742912

# Row 5 (SIZE)
$ws.Range("D5").Value = 0.1370651265057656
$ws.Range("E5").Value = 0.003006012024048044

# Row 6 (USMV)
$ws.Range("D6").Value = 0.131191493474018
$ws.Range("E6").Value = 0.008999704927707297

# Row 7 (VLUE)
$ws.Range("D7").Value = 0.1461154157883357
$ws.Range("E7").Value = 0.00454765360425724

# Row 8 (SNSR)
$ws.Range("D8").Value = 0.1279010606191076
$ws.Range("E8").Value = 0.009382566585956376

# Row 9 (FITE)
$ws.Range("D9").Value = 0.1301941876804396
$ws.Range("E9").Value = -0.0009311293180350733

# Row 10 (Total) - only Percent Change changes, Weight stays 1
$ws.Range("E10").Value = 0.008282588378830003
